$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 15 updates
$ws.Range("A15").Value = 112322598
$ws.Range("P15").Value = "Gräsviggen, Vrm"
$ws.Range("Q15").Value = 376602
$ws.Range("R15").Value = 6700269
$ws.Range("AW15").Value = "Helena Malmestrand"
$ws.Range("AX15").Value = "Helena Malmestrand"

# Row 16 updates
$ws.Range("A16").Value = 112322551
$ws.Range("B16").Value = 77671
$ws.Range("E16").Value = 185
$ws.Range("F16").Value = "Violettgrå tagellav"
$ws.Range("G16").Value = "Bryoria nadvornikiana"
$ws.Range("H16").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q16").Value = 376720
$ws.Range("R16").Value = 6700740

# Row 17 updates
$ws.Range("A17").Value = 112323374
$ws.Range("B17").Value = 77636
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = "Garnlav"
$ws.Range("G17").Value = "Alectoria sarmentosa"
$ws.Range("H17").Value = "(Ach.) Ach."
$ws.Range("P17").Value = "Jonsmyren, Vrm"
$ws.Range("Q17").Value = 376704
$ws.Range("R17").Value = 6700726
$ws.Range("AW17").Value = "anders tedeholm"
$ws.Range("AX17").Value = "anders tedeholm"

# Row 36 updates
$ws.Range("A36").Value = 112323387
$ws.Range("B36").Value = 77671
$ws.Range("E36").Value = 185
$ws.Range("F36").Value = "Violettgrå tagellav"
$ws.Range("G36").Value = "Bryoria nadvornikiana"
$ws.Range("H36").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("P36").Value = "Jonsmyren, Vrm"
$ws.Range("Q36").Value = 376577
$ws.Range("R36").Value = 6700273
$ws.Range("AW36").Value = "anders tedeholm"
$ws.Range("AX36").Value = "anders tedeholm"

# Row 37 updates
$ws.Range("A37").Value = 112323371
$ws.Range("B37").Value = 77636
$ws.Range("E37").Value = 6425
$ws.Range("F37").Value = "Garnlav"
$ws.Range("G37").Value = "Alectoria sarmentosa"
$ws.Range("H37").Value = "(Ach.) Ach."
$ws.Range("Q37").Value = 376715
$ws.Range("R37").Value = 6700826

# Row 38 updates
$ws.Range("A38").Value = 112322559
$ws.Range("P38").Value = "Gräsviggen, Vrm"
$ws.Range("Q38").Value = 376483
$ws.Range("R38").Value = 6700251
$ws.Range("AW38").Value = "Helena Malmestrand"
$ws.Range("AX38").Value = "Helena Malmestrand"
